$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.173.25"
Set-TextValue $ws.Range("E2") "  +3.56%  "
Set-TextValue $ws.Range("D3") "1.895.08"
Set-TextValue $ws.Range("E3") "  +3.98%  "
Set-TextValue $ws.Range("D4") "0.9988"
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "246.59"
Set-TextValue $ws.Range("E5") "  +0.15%  "
Set-TextValue $ws.Range("D6") "0.9984"
Set-TextValue $ws.Range("E6") "  -0.24%  "
Set-TextValue $ws.Range("D7") "0.4987"
Set-TextValue $ws.Range("E7") "  +1.19%  "
Set-TextValue $ws.Range("D8") "44.82"
Set-TextValue $ws.Range("E8") "  +0.88%  "
Set-TextValue $ws.Range("D9") "0.2941"
Set-TextValue $ws.Range("E9") "  +5.90%  "
Set-TextValue $ws.Range("D10") "0.06656"
Set-TextValue $ws.Range("E10") "  +4.30%  "
Set-TextValue $ws.Range("D11") "1.892.23"
Set-TextValue $ws.Range("E11") "  +3.83%  "
Set-TextValue $ws.Range("D12") "16.98"
Set-TextValue $ws.Range("E12") "  +1.95%  "
Set-TextValue $ws.Range("E13") "  +1.23%  "
Set-TextValue $ws.Range("D14") "0.6762"
Set-TextValue $ws.Range("E14") "  +4.86%  "
Set-TextValue $ws.Range("D15") "86.08"
Set-TextValue $ws.Range("E15") "  +2.41%  "
Set-TextValue $ws.Range("D16") "4.850"
Set-TextValue $ws.Range("E16") "  +3.17%  "
Set-TextValue $ws.Range("D17") "30.158.90"
Set-TextValue $ws.Range("E17") "  +3.52%  "
Set-TextValue $ws.Range("D18") "0.000007989"
Set-TextValue $ws.Range("E18") "  +9.31%  "
Set-TextValue $ws.Range("D19") "0.9974"
Set-TextValue $ws.Range("E19") "  -0.21%  "
Set-TextValue $ws.Range("D20") "12.86"
Set-TextValue $ws.Range("E20") "  +5.19%  "
Set-TextValue $ws.Range("D21") "2.136.08"
Set-TextValue $ws.Range("E21") "  +3.85%  "
Set-TextValue $ws.Range("D22") "0.9991"
Set-TextValue $ws.Range("E22") "  -0.08%  "
Set-TextValue $ws.Range("D23") "4.780"
Set-TextValue $ws.Range("E23") "  +5.17%  "
Set-TextValue $ws.Range("D24") "5.641"
Set-TextValue $ws.Range("E24") "  +5.07%  "
Set-TextValue $ws.Range("D25") "9.138"
Set-TextValue $ws.Range("E25") "  +3.35%  "
Set-TextValue $ws.Range("D26") "147.84"
Set-TextValue $ws.Range("E26") "  +2.93%  "
Set-TextValue $ws.Range("D27") "134.37"
Set-TextValue $ws.Range("E27") "  +3.05%  "
Set-TextValue $ws.Range("D28") "16.82"
Set-TextValue $ws.Range("E28") "  +2.58%  "
Set-TextValue $ws.Range("D29") "1.944"
Set-TextValue $ws.Range("E29") "  +3.17%  "
Set-TextValue $ws.Range("D30") "1.378"
Set-TextValue $ws.Range("E30") "  -1.59%  "
Set-TextValue $ws.Range("D31") "4.211"
Set-TextValue $ws.Range("E31") "  +2.22%  "
Set-TextValue $ws.Range("D32") "0.08734"
Set-TextValue $ws.Range("E32") "  +4.87%  "
Set-TextValue $ws.Range("D33") "3.965"
Set-TextValue $ws.Range("E33") "  +5.02%  "
Set-TextValue $ws.Range("D34") "0.05152"
Set-TextValue $ws.Range("E34") "  +4.54%  "
Set-TextValue $ws.Range("D35") "1.124"
Set-TextValue $ws.Range("E35") "  +2.64%  "
Set-TextValue $ws.Range("D36") "0.7075"
Set-TextValue $ws.Range("E36") "  +5.82%  "
Set-TextValue $ws.Range("D37") "2.666"
Set-TextValue $ws.Range("E37") "  -1.10%  "
Set-TextValue $ws.Range("D38") "2.777"
Set-TextValue $ws.Range("E38") "  +3.62%  "
Set-TextValue $ws.Range("D39") "2.236"
Set-TextValue $ws.Range("E39") "  -2.30%  "
Set-TextValue $ws.Range("D40") "0.9417"
Set-TextValue $ws.Range("E40") "  -0.89%  "
Set-TextValue $ws.Range("D41") "0.01662"
Set-TextValue $ws.Range("E41") "  +4.55%  "
Set-TextValue $ws.Range("D42") "6.090"
Set-TextValue $ws.Range("E42") "  -1.01%  "
Set-TextValue $ws.Range("D45") "102.88"
Set-TextValue $ws.Range("E45") "  +0.96%  "
Set-TextValue $ws.Range("D46") "7.525"
Set-TextValue $ws.Range("E46") "  +5.19%  "
Set-TextValue $ws.Range("D47") "0.1262"
Set-TextValue $ws.Range("E47") "  +3.71%  "
Set-TextValue $ws.Range("D48") "0.05728"
Set-TextValue $ws.Range("E48") "  +3.25%  "
Set-TextValue $ws.Range("D49") "32.76"
Set-TextValue $ws.Range("E49") "  +3.24%  "
Set-TextValue $ws.Range("D50") "8.290"
Set-TextValue $ws.Range("E50") "  +2.64%  "
Set-TextValue $ws.Range("E51") "  +4.08%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "0.9971"
Set-TextValue $ws.Range("E43") "  -0.32%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D44") "0.4213"
Set-TextValue $ws.Range("E44") "  +3.78%  "
